$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 130 (pushes existing rows 130-247 down to 131-248).
$ws.Rows.Item(130).Insert()

# Populate the newly inserted row 130 with the new weekly price record.
$ws.Cells.Item(130, 1).Value  = 11
$ws.Cells.Item(130, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(130, 3).Value  = "Bíobío"
$ws.Cells.Item(130, 4).Value  = 44566
$ws.Cells.Item(130, 5).Value  = 8
$ws.Cells.Item(130, 6).Value  = 100114014
$ws.Cells.Item(130, 7).Value  = "Betarraga"
$ws.Cells.Item(130, 8).Value  = "Sin especificar"
$ws.Cells.Item(130, 9).Value  = "Primera"
$ws.Cells.Item(130, 10).Value = 500
$ws.Cells.Item(130, 11).Value = 600
$ws.Cells.Item(130, 12).Value = 650
$ws.Cells.Item(130, 13).Value = 620
$ws.Cells.Item(130, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(130, 15).Value = "Región Metropolitana"
$ws.Cells.Item(130, 16).Value = 124
$ws.Cells.Item(130, 17).Value = 5
$ws.Cells.Item(130, 18).Value = "Hortaliza"
